# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.158.21'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').Value = '1.654.36'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.30'
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5114'
$ws.Range('E6').Value = '  -3.43%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06403'
$ws.Range('E9').Value = '  -3.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.93'
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07821'
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.282'
$ws.Range('E12').Value = '  -5.27%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.882.14'
$ws.Range('E13').Value = '  -3.30%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.632.60'
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5525'
$ws.Range('E15').Value = '  -5.30%  '
$ws.Range('D16').Value = '0.0₅8023'
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.88'
$ws.Range('E17').Value = '  -6.19%  '
$ws.Range('D18').Value = '26.178.87'
$ws.Range('E18').Value = '  -4.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '209.05'
$ws.Range('E20').Value = '  -6.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.417'
$ws.Range('E21').Value = '  -4.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.06'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.027'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.06'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.727'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1166'
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.971'
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05117'
$ws.Range('E30').Value = '  -4.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.243'
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.353'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.221'
$ws.Range('E33').Value = '  -6.32%  '
$ws.Range('E34').Value = '  -5.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.748'
$ws.Range('E35').Value = '  -4.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.373'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9285'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5680'
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').Value = '1.152.65'
$ws.Range('E39').Value = '  +5.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01590'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8335'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.642'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.07'
$ws.Range('E44').Value = '  -1.17%  '
$ws.Range('D45').Value = '1.792.62'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4546'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.63'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.842'
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05035'
$ws.Range('E51').Value = '  -3.93%  '
